# (#230) corrected test data
# Fix two erroneous parameter-string values on the "cancertype" sheet:
#  - L3 (ParamCancerTypeStageSideEffects) used "st=" (SubType) instead of "stg=" (Stage)
#  - I4 (ParamCancerTypeSideEffects) incorrectly included "st=" and "stg=" params
#    that don't belong in a SideEffects-only param string

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cancertype")

$ws.Range("I4").Value = "fin=C28306&loc=0&rl=2&t=C9291"
$ws.Range("L3").Value = "fin=C115197&loc=0&rl=2&stg=C94774&t=C4872"

# Leave the selection on the corrected cell, matching the saved view state
$ws.Range("L3").Select()
